$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the label in A53 from "min" to "min/max"
$ws.Range("A53").Value = "min/max"

# Set the new formulas in B53 and C53
$ws.Range("B53").Formula = "=MAX(B1:B51)"
$ws.Range("C53").Formula = "=MIN(C1:C51)"

# Adjust column widths
$ws.Columns.Item(2).ColumnWidth = 23.5703125
$ws.Columns.Item(3).ColumnWidth = 29.140625

# Update the view: scroll position and selection
$ws.Range("C53").Select()
$excel.ActiveWindow.ScrollRow = 43
$excel.ActiveWindow.ScrollColumn = 1
